# "DMP real time added"
# Adds a new "DMPRealtimecontainers" column (L) to the PegaTestData sheet
# with the value "CSM", matching the look of the existing header columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header (row 1) and data value (row 2) in column L
$ws.Range("L1").Value = "DMPRealtimecontainers"
$ws.Range("L2").Value = "CSM"

# Give the new header cell the same (green) fill as the other header cells
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Size the new column similarly to its neighbours
$ws.Columns("L").ColumnWidth = 21.6

# Leave the selection where the workbook was last left
$ws.Range("I8").Select() | Out-Null
